$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F6").Value = 11337
    $ws.Range("F8").Value = 104
    $ws.Range("F18").Value = 321
    $ws.Range("F19").Value = 1259
    $ws.Range("F20").Value = 65
}

$wb.Save()
